$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "69.856.24"
$ws.Range("E2").Value = "  -1.85%  "
Set-TextCell "D3" "3.756.30"
$ws.Range("E3").Value = "  +2.96%  "
Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextCell "D5" "621.18"
$ws.Range("E5").Value = "  +3.89%  "
Set-TextCell "D6" "180.45"
$ws.Range("E6").Value = "  -1.06%  "
Set-TextCell "D7" "3.755.32"
$ws.Range("E7").Value = "  +3.09%  "
Set-TextCell "D8" "0.999"
$ws.Range("E8").Value = "  +0.09%  "
Set-TextCell "D9" "0.533"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E11").Value = "  -5.55%  "
$ws.Range("E12").Value = "  -1.94%  "
Set-TextCell "D13" "41.44"
$ws.Range("E13").Value = "  +1.19%  "
Set-TextCell "D14" "0.0000257"
$ws.Range("E14").Value = "  +0.92%  "
Set-TextCell "D15" "4.375.97"
$ws.Range("E15").Value = "  +2.75%  "
Set-TextCell "D16" "3.756.03"
$ws.Range("E16").Value = "  +3.04%  "
Set-TextCell "D17" "69.907.14"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +0.78%  "
Set-TextCell "D20" "16.61"
$ws.Range("E20").Value = "  -2.32%  "
Set-TextCell "D21" "508.19"
$ws.Range("E21").Value = "  -1.55%  "
Set-TextCell "D22" "9.58"
$ws.Range("E22").Value = "  +4.58%  "
Set-TextCell "D23" "0.724"
$ws.Range("E23").Value = "  -2.67%  "
Set-TextCell "D24" "2.53"
$ws.Range("E24").Value = "  +1.16%  "
Set-TextCell "D25" "86.98"
$ws.Range("E25").Value = "  -0.89%  "
Set-TextCell "D26" "13.10"
$ws.Range("E26").Value = "  -3.63%  "
Set-TextCell "D27" "11.09"
$ws.Range("E27").Value = "  +0.54%  "
Set-TextCell "D28" "0.0000134"
$ws.Range("E28").Value = "  +21.12%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -1.28%  "
Set-TextCell "D31" "2.93"
$ws.Range("E31").Value = "  +5.41%  "
Set-TextCell "D32" "7.90"
$ws.Range("E32").Value = "  -3.13%  "
Set-TextCell "D33" "30.96"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("E34").Value = "  -1.51%  "
Set-TextCell "D35" "1.00"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +5.18%  "
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("E40").Value = "  -2.27%  "
Set-TextCell "D41" "50.02"
$ws.Range("E41").Value = "  -2.17%  "
Set-TextCell "D42" "45.49"
$ws.Range("E42").Value = "  +0.73%  "
Set-TextCell "D43" "425.04"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("E44").Value = "  -1.41%  "
Set-TextCell "D45" "3.017.07"
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -1.92%  "
Set-TextCell "D48" "27.49"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("E49").Value = "  -0.06%  "
Set-TextCell "D50" "137.59"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  -0.01%  "
